$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsSummary = $wb.Worksheets.Item("Summary")

$str53 = "[15, 16, 9, 10, 109, 0.67601085, 0.513422, 0.2635864, 0.21835516, 0.268843, 0.41236666, 0.65258217, 0.6229082, 0.29250848, 0.528681, 0.75332034, 0.8282525, 0.5150838, 0.9204639, 0.80069244, 0.12044735, 0.2903722, 0.71965295, 0.22852382, 0.7829796, 0.30210385, 0.5443946, 0.58138466, 0.5813501, 0.34919834, 0.20638384, 0.5830844, 0.50312436, 0.54213613, 0.50839067, 0.7508806, 0.41311, 0.71007115, 0.40801752, 0.80320704, 0.9690926, 0.77518743, 0.56275654, 0.4808577, 0.80989414, 0.049788415, 0.4248149, 0.2715052, 0.13420284, 0.27404973, 0.79519796, 0.9279251]"
$str54 = "[10, 8, 13, 7, 125, 0.56386465, 0.4863702, 0.61907876, 0.61731106, 0.5255747, 0.66833234, 0.5591835, 0.49343428, 0.50020045, 0.36494184, 0.45790708, 0.8200666, 0.81914735, 0.7505008, 0.7158586, 0.014160433, 0.2507621, 0.6299417, 0.43736744, 0.8989186, 0.29443944, 0.9440727, 0.09307182, 0.1940778, 0.3310041, 0.35835853, 0.92452115, 0.33856618, 0.4840274, 0.55174744, 0.41118985, 0.6761216, 0.9114944, 0.5781735, 0.35986722, 0.80815345, 0.6171422, 0.95544565, 0.0035003424, 0.84136665, 0.040762722, 0.57797354, 0.5832304, 0.4854872, 0.06584303, 0.29277968, 0.8486911]"

# Data sheet row 9
$wsData.Cells.Item(9, 1).Value = 6.0
$wsData.Cells.Item(9, 2).Value = 815.86664
$wsData.Cells.Item(9, 3).Value = 278.2489
$wsData.Cells.Item(9, 4).Value = 170.79956
$wsData.Cells.Item(9, 5).Value = 815.86664
$wsData.Cells.Item(9, 6).Value = 130.44612
$wsData.Cells.Item(9, 7).Value = 30.0
$wsData.Cells.Item(9, 8).Value = $str53
$wsData.Cells.Item(9, 9).Value = 94.26667
$wsData.Cells.Item(9, 10).Value = 32.628456
$wsData.Cells.Item(9, 11).Value = 418.66666
$wsData.Cells.Item(9, 12).Value = 21.745129
$wsData.Cells.Item(9, 13).Value = 163.93333
$wsData.Cells.Item(9, 14).Value = 60.595284
$wsData.Cells.Item(9, 15).Value = 0.0
$wsData.Cells.Item(9, 16).Value = 0.0
$wsData.Cells.Item(9, 17).Value = 32.033333
$wsData.Cells.Item(9, 18).Value = 7.5907536
$wsData.Cells.Item(9, 19).Value = 32.033333
$wsData.Cells.Item(9, 20).Value = 7.5907536
$wsData.Cells.Item(9, 21).Value = 13.7
$wsData.Cells.Item(9, 22).Value = 6.6027684
$wsData.Cells.Item(9, 23).Value = 0.0
$wsData.Cells.Item(9, 24).Value = 0.0
$wsData.Cells.Item(9, 25).Value = 32.033333
$wsData.Cells.Item(9, 26).Value = 7.5907536
$wsData.Cells.Item(9, 27).Value = 0.33333334
$wsData.Cells.Item(9, 28).Value = 1.2954384
$wsData.Cells.Item(9, 29).Value = 28.866667
$wsData.Cells.Item(9, 30).Value = 27.660172

# Data sheet row 10
$wsData.Cells.Item(10, 1).Value = 7.0
$wsData.Cells.Item(10, 2).Value = 920.86664
$wsData.Cells.Item(10, 3).Value = 323.1889
$wsData.Cells.Item(10, 4).Value = 193.27455
$wsData.Cells.Item(10, 5).Value = 920.86664
$wsData.Cells.Item(10, 6).Value = 136.95172
$wsData.Cells.Item(10, 7).Value = 30.0
$wsData.Cells.Item(10, 8).Value = $str54
$wsData.Cells.Item(10, 9).Value = 41.766666
$wsData.Cells.Item(10, 10).Value = 20.5723
$wsData.Cells.Item(10, 11).Value = 407.96667
$wsData.Cells.Item(10, 12).Value = 30.68582
$wsData.Cells.Item(10, 13).Value = 240.06667
$wsData.Cells.Item(10, 14).Value = 74.35606
$wsData.Cells.Item(10, 15).Value = 0.033333335
$wsData.Cells.Item(10, 16).Value = 0.18257418
$wsData.Cells.Item(10, 17).Value = 43.166668
$wsData.Cells.Item(10, 18).Value = 7.926748
$wsData.Cells.Item(10, 19).Value = 43.166668
$wsData.Cells.Item(10, 20).Value = 7.926748
$wsData.Cells.Item(10, 21).Value = 46.4
$wsData.Cells.Item(10, 22).Value = 18.225788
$wsData.Cells.Item(10, 23).Value = 0.0
$wsData.Cells.Item(10, 24).Value = 0.0
$wsData.Cells.Item(10, 25).Value = 43.166668
$wsData.Cells.Item(10, 26).Value = 7.926748
$wsData.Cells.Item(10, 27).Value = 0.3
$wsData.Cells.Item(10, 28).Value = 0.83666
$wsData.Cells.Item(10, 29).Value = 54.833332
$wsData.Cells.Item(10, 30).Value = 23.64549

# Data sheet row 11
$wsData.Cells.Item(11, 1).Value = 8.0
$wsData.Cells.Item(11, 2).Value = 920.86664
$wsData.Cells.Item(11, 3).Value = 369.25333
$wsData.Cells.Item(11, 4).Value = 149.78954
$wsData.Cells.Item(11, 5).Value = 920.86664
$wsData.Cells.Item(11, 6).Value = 136.95172
$wsData.Cells.Item(11, 7).Value = 30.0
$wsData.Cells.Item(11, 8).Value = $str54
$wsData.Cells.Item(11, 9).Value = 41.766666
$wsData.Cells.Item(11, 10).Value = 20.5723
$wsData.Cells.Item(11, 11).Value = 407.96667
$wsData.Cells.Item(11, 12).Value = 30.68582
$wsData.Cells.Item(11, 13).Value = 240.06667
$wsData.Cells.Item(11, 14).Value = 74.35606
$wsData.Cells.Item(11, 15).Value = 0.033333335
$wsData.Cells.Item(11, 16).Value = 0.18257418
$wsData.Cells.Item(11, 17).Value = 43.166668
$wsData.Cells.Item(11, 18).Value = 7.926748
$wsData.Cells.Item(11, 19).Value = 43.166668
$wsData.Cells.Item(11, 20).Value = 7.926748
$wsData.Cells.Item(11, 21).Value = 46.4
$wsData.Cells.Item(11, 22).Value = 18.225788
$wsData.Cells.Item(11, 23).Value = 0.0
$wsData.Cells.Item(11, 24).Value = 0.0
$wsData.Cells.Item(11, 25).Value = 43.166668
$wsData.Cells.Item(11, 26).Value = 7.926748
$wsData.Cells.Item(11, 27).Value = 0.3
$wsData.Cells.Item(11, 28).Value = 0.83666
$wsData.Cells.Item(11, 29).Value = 54.833332
$wsData.Cells.Item(11, 30).Value = 23.64549

# Summary sheet row 8
$wsSummary.Cells.Item(8, 1).Value = 6.0
$wsSummary.Cells.Item(8, 2).Value = 815.86664
$wsSummary.Cells.Item(8, 3).Value = 278.2489

# Summary sheet row 9
$wsSummary.Cells.Item(9, 1).Value = 7.0
$wsSummary.Cells.Item(9, 2).Value = 920.86664
$wsSummary.Cells.Item(9, 3).Value = 323.1889

# Summary sheet row 10
$wsSummary.Cells.Item(10, 1).Value = 8.0
$wsSummary.Cells.Item(10, 2).Value = 920.86664
$wsSummary.Cells.Item(10, 3).Value = 369.25333

